$wb = $excel.ActiveWorkbook

# --- Update Metadata sheet: Date and Count values ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-09-13T20:57:31+00:00"
$meta.Range("B22").Value = "3"

# --- Add new concept row to the Concepts sheet ---
$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Range("A4").Value = "1"
$concepts.Range("B4").Value = "unknown"
$concepts.Range("C4").Value = "Unknown"
$concepts.Range("A4:D4").Style = $concepts.Range("A2:D2").Style
